# Update "想去人数" (interested-count) figures on the 展览 (Exhibitions),
# 演出 (Performances) and 全部类型 (All types) sheets, as published by the
# latest scrape run (gh-pages output regenerated at 456a3b4).
#
# 本地生活 (sheet 3) has no data rows in this run, so it needs no edits.

$wb = $excel.ActiveWorkbook

# --- 展览 ("展览" sheet) ---------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value  = 561    # 苏州·归离之缘原神only展
$wsExpo.Range("F4").Value  = 1132   # 张家港·幻想物语动漫游戏嘉年华02
$wsExpo.Range("F10").Value = 16303  # 苏州·萤火国潮文化节动漫品牌博览会
$wsExpo.Range("F14").Value = 6355   # 苏州·OCG国潮动漫嘉年华
$wsExpo.Range("F15").Value = 636    # 苏州·白日梦想7.20全职猎人ONLY展
$wsExpo.Range("F19").Value = 124    # 苏州·AME动漫嘉年华
$wsExpo.Range("F21").Value = 32     # 苏州·第五人格only·盛典
$wsExpo.Range("F23").Value = 638    # 张家港·喵言动漫游戏嘉年华
$wsExpo.Range("F24").Value = 31     # 昆山·第七届·xcy新次元动漫嘉年华-狂欢盛典
$wsExpo.Range("F27").Value = 221    # 常熟·ACG动漫游戏嘉年华
$wsExpo.Range("F28").Value = 892    # 常熟·CDW.动漫展05
$wsExpo.Range("F29").Value = 53     # 苏州·代号鸢only茶话会-星渡咖啡
$wsExpo.Range("F32").Value = 11313  # 苏州·ICAN summer World动漫品牌夏游节
$wsExpo.Range("F33").Value = 1245   # 苏州·第二届Redamancy动漫游戏嘉年华
$wsExpo.Range("F35").Value = 148    # 苏州·Good jump ACG中秋嘉年华动漫国潮文化节
$wsExpo.Range("F36").Value = 205    # 苏州·I COME ACG动漫品牌博览会

# --- 演出 ("演出" sheet) ---------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 21      # 苏州·爱乐之城·经典电影作品音乐会

# --- 全部类型 ("全部类型" sheet) --------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 561     # 苏州·归离之缘原神only展
$wsAll.Range("F4").Value  = 1132    # 张家港·幻想物语动漫游戏嘉年华02
$wsAll.Range("F10").Value = 16303   # 苏州·萤火国潮文化节动漫品牌博览会
$wsAll.Range("F14").Value = 6355    # 苏州·OCG国潮动漫嘉年华
$wsAll.Range("F15").Value = 636     # 苏州·白日梦想7.20全职猎人ONLY展
$wsAll.Range("F19").Value = 124     # 苏州·AME动漫嘉年华
$wsAll.Range("F21").Value = 32      # 苏州·第五人格only·盛典
$wsAll.Range("F23").Value = 638     # 张家港·喵言动漫游戏嘉年华
$wsAll.Range("F24").Value = 31      # 昆山·第七届·xcy新次元动漫嘉年华-狂欢盛典
$wsAll.Range("F27").Value = 221     # 常熟·ACG动漫游戏嘉年华
$wsAll.Range("F28").Value = 892     # 常熟·CDW.动漫展05
$wsAll.Range("F29").Value = 53      # 苏州·代号鸢only茶话会-星渡咖啡
$wsAll.Range("F32").Value = 21      # 苏州·爱乐之城·经典电影作品音乐会
$wsAll.Range("F33").Value = 11313   # 苏州·ICAN summer World动漫品牌夏游节
$wsAll.Range("F34").Value = 1245    # 苏州·第二届Redamancy动漫游戏嘉年华
$wsAll.Range("F36").Value = 148     # 苏州·Good jump ACG中秋嘉年华动漫国潮文化节
$wsAll.Range("F37").Value = 205     # 苏州·I COME ACG动漫品牌博览会
